# The commit adds one new weekly price-report row for
# "Comercializadora del Agro de Limarí - Arveja Verde" on top of the
# existing daily log. The new observation is inserted as row 9 (the rows
# are apparently kept in reverse-chronological/insert order), pushing the
# previous rows 9..88 down to 10..89 and growing the used range from
# A1:R88 to A1:R89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts existing rows 9-88 down to
# 10-89 (cell formatting/styles carry down with them automatically) and
# extends the sheet's used range accordingly.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new observation. All the
# "constant" descriptive columns reuse the same values found throughout
# the rest of the table.
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 45092
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Perfection"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 27000
$ws.Range("L9").Value = 29000
$ws.Range("M9").Value = 28000
$ws.Range("N9").Value = "$/malla 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 1120
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
